# Defect Log.xlsx update
# Nhom ThanhChV, LinhTA, HuyDV, DucNH cap nhat fie Defect Log.xlsx sau khi sua loi.
#
# For defect rows 15-26: mark Status as "Corrected" and fill in
# Fixed Date / Closed Date with 2011-10-22 (serial 40838), matching the
# date already used for the other completed defects in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixedDate = 40838

foreach ($r in 15..26) {
    $ws.Range("E$r").Value = "Corrected"

    $ws.Range("S$r").Value = $fixedDate

    # Match the "Closed Date" cell's formatting/border to the "Fixed Date"
    # cell in the same row (copy format only), then set its value.
    $ws.Range("S$r").Copy() | Out-Null
    $ws.Range("T$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("T$r").Value = $fixedDate
}

$excel.CutCopyMode = 0

# Restore the view/selection state left by the editors after finishing
# the update.
$ws.Range("E16").Select() | Out-Null
